$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs","Bmp6","Bmpr1b","FAPs",3,1,8.929813666666668,26.789441,0.3579859341865942,0.3579859341865942,3,1,3.510190333333334,10.530571,0.7927501397588634,0.7927501397588635,31.34534561120123,282.1081105008111,0.2837933993581298,0.2837933993581299),
  @("ECs","Bmp6","Bmpr1b","MuSCs",3,1,8.929813666666668,26.789441,0.3579859341865942,0.3579859341865942,3,1,0.8869683333333332,2.660905,0.2003151406163121,0.2003151406163121,7.920461944900556,71.284157504105,0.07171000274524945,0.07171000274524947),
  @("ECs","Bmp6","Bmpr1b","Resolving-Mac",3,1,8.929813666666668,26.789441,0.3579859341865942,0.3579859341865942,1,0.3333333333333333,0.030706,0.092118,0.006934719624824425,0.006934719624824427,0.2741988584486667,2.467789726038001,0.00248253208321488,0.002482532083214881),
  @("FAPs","Bmp6","Bmpr1b","FAPs",3,1,0.7684289999999999,2.305287,0.03080543264277933,0.03080543264277933,3,1,3.510190333333334,10.530571,0.7927501397588634,0.7927501397588635,2.697332047653,24.275988428877,0.02442101103289557,0.02442101103289557),
  @("FAPs","Bmp6","Bmpr1b","MuSCs",3,1,0.7684289999999999,2.305287,0.03080543264277933,0.03080543264277933,3,1,0.8869683333333332,2.660905,0.2003151406163121,0.2003151406163121,0.6815721894149999,6.134149704734999,0.006170794571584673,0.006170794571584674),
  @("FAPs","Bmp6","Bmpr1b","Resolving-Mac",3,1,0.7684289999999999,2.305287,0.03080543264277933,0.03080543264277933,1,0.3333333333333333,0.030706,0.092118,0.006934719624824425,0.006934719624824427,0.023595380874,0.212358427866,0.0002136270382990888,0.0002136270382990889),
  @("MuSCs","Bmp6","Bmpr1b","FAPs",3,1,15.246351,45.739053,0.6112086331706265,0.6112086331706265,3,1,3.510190333333334,10.530571,0.7927501397588634,0.7927501397588635,53.517593898807,481.658345089263,0.484535729367838,0.4845357293678381),
  @("MuSCs","Bmp6","Bmpr1b","MuSCs",3,1,15.246351,45.739053,0.6112086331706265,0.6112086331706265,3,1,0.8869683333333332,2.660905,0.2003151406163121,0.2003151406163121,13.523030535885,121.707274822965,0.122434343299478,0.122434343299478),
  @("MuSCs","Bmp6","Bmpr1b","Resolving-Mac",3,1,15.246351,45.739053,0.6112086331706265,0.6112086331706265,1,0.3333333333333333,0.030706,0.092118,0.006934719624824425,0.006934719624824427,0.468154453806,4.213390084254,0.004238560503310456,0.004238560503310457)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
